$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the four new data rows (38-41) with updated fastq/metadata values.
#    Columns A and D hold "MM.DD.YY"-looking text that Excel's COM layer
#    would otherwise auto-convert to a date serial, so each of those cells
#    is pre-formatted as Text ("@") before the value is written, and the
#    formatting is reset back to Normal afterwards so no stray number
#    format lingers on the cell.
#    Row 40 is the only row that introduces brand-new shared strings
#    ("11.01.18" / "11.02.18"); column D is written before column A there
#    so the new strings land in the same table order as the source edit.
# ---------------------------------------------------------------------------

# Row 38
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "08.13.18"
$ws.Range("B38").Value = "H.BROWN"
$ws.Range("C38").Value = 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "08.14.18"
$ws.Range("E38").Value = "H.BROWN"
$ws.Range("F38").Value = 38
$ws.Range("G38").Value = "E7420L"

# Row 39
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "10.16.18"
$ws.Range("B39").Value = "H.BROWN"
$ws.Range("C39").Value = 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.18.18"
$ws.Range("E39").Value = "H.BROWN"
$ws.Range("F39").Value = 39
$ws.Range("G39").Value = "E7420L"

# Row 40 (introduces the two brand new shared strings; D before A)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.02.18"
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "11.01.18"
$ws.Range("B40").Value = "H.BROWN"
$ws.Range("C40").Value = 40
$ws.Range("E40").Value = "H.BROWN"
$ws.Range("F40").Value = 40
$ws.Range("G40").Value = "E7420L"

# Row 41
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "10.16.18"
$ws.Range("B41").Value = "H.BROWN"
$ws.Range("C41").Value = 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.18.18"
$ws.Range("E41").Value = "H.BROWN"
$ws.Range("F41").Value = 41
$ws.Range("G41").Value = "E7420L"

# ---------------------------------------------------------------------------
# 2. Strip the bold/bordered header style from row 1 (A1:F1) and the
#    wrap-text style from the whole "s1cDNAProtocol" column (G2:G41), so
#    every cell reverts to the default/Normal style (no "s" attribute).
#    This also clears the temporary Text ("@") format applied above to the
#    new rows' A/D date-like cells.
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").Style = "Normal"
$ws.Range("A38:G41").Style = "Normal"
$ws.Range("G2:G41").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Drop the custom row height on the pre-existing data rows so they fall
#    back to the sheet's default row height, matching rows 38-41 which were
#    never given an explicit height.
# ---------------------------------------------------------------------------
$ws.Range("A2:A37").EntireRow.AutoFit()

# ---------------------------------------------------------------------------
# 4. Update the selection to match the post-edit state (a "select all"
#    style range covering columns A:I).
# ---------------------------------------------------------------------------
$ws.Range("A1:I1048576").Select()
